$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10 (shifts old rows 10-14 down to 11-15)
$ws.Rows.Item(10).Insert()

# Populate the new row 10 values (string order matters for shared-string table: 1N4001 then Diodes (x4))
$ws.Cells.Item(10, 2).Value = "1N4001"
$ws.Cells.Item(10, 1).Value = "Diodes (x4)"
$ws.Cells.Item(10, 3).Formula = "=0.21/10*4"
$ws.Cells.Item(10, 4).Value = "Futurlec"
$ws.Cells.Item(10, 5).Value = "Ok"

# Copy style from the cost cell above (row 9) to maintain currency formatting
$ws.Cells.Item(9, 3).Copy()
$ws.Cells.Item(10, 3).PasteSpecial(-4122)  # xlPasteFormats

# Resize the table to include the new row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E15"))

$excel.CalculateFullRebuild()
